# "Finished Design thinking section" — update the design-thinking summary
# cell (I33, merged I33:L33) with its final wording, let the text wrap over
# two lines, and resize the row/columns around it so the sheet reads cleanly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. New wording for the merged header cell I33:L33.
$ws.Range("I33").Value = "Average total time when feedback`nwas first or second"

# 2. Let the new (two-line) text wrap, and grow the row to fit it.
$ws.Range("I33").WrapText = $true
$ws.Rows.Item(33).RowHeight = 34.5

# 3. The J/K/L summary columns no longer need to be as wide now that the
#    header above them wraps instead of stretching the columns.
$ws.Columns.Item(10).ColumnWidth = 5.166666666666667
$ws.Columns.Item(11).ColumnWidth = 4.666666666666667
$ws.Columns.Item(12).ColumnWidth = 7.833333333333333

# 4. Row 46's accuracy average was pointing at the wrong block of rows
#    (A146:A151 overlaps the "feedback second" set already used in J46) -
#    fix it to reference A152:A157, matching the other L-column formulas.
$ws.Range("L46").Formula = "=AVERAGE(A8:A13,A40:A45,A72:A77,A152:A157)"

# 5. Leave the cursor where the author left it.
$ws.Range("O33").Select() | Out-Null
